# Regenerate s_vals data to filter save games.
# Updates columns B,C,D,E,G for rows 2-9. Column F (Win) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.01514828764759746; C = 0.3127903958511391;  D = 0.1575252929769615;  E = 8.660232485948974;  G = 9.145696462424672  }
    3 = @{ B = 0.3048080303191223;  C = 0.04240448674262143; D = 3.900430680208489;   E = 0.496779210170732;  G = 4.744422407440965  }
    4 = @{ B = 0.127881588408715;   C = 0.3127903958511391;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 1.094976487407548  }
    5 = @{ B = 3.230985683306322;   C = 1.667794583268128;   D = 0.8054896365839992;  E = 0.496779210170732;  G = 6.201049113329182  }
    6 = @{ B = 3.230985683306322;   C = 1.667794583268128;   D = 0.8054896365839992;  E = 8.660232485948974;  G = 14.36450238910742  }
    7 = @{ B = 0.0001488876196638067; C = 0.04240448674262143; D = 3.900430680208489; E = 8.660232485948974;  G = 12.60321654051975  }
    8 = @{ B = 0.127881588408715;   C = 1.667794583268128;   D = 0.1575252929769615;  E = 8.660232485948974;  G = 10.61343395060278  }
    9 = @{ B = 0.01514828764759746; C = 1.667794583268128;   D = 0.1575252929769615;  E = 0.496779210170732;  G = 2.337247374063419  }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
